# Append three new gene-expression rows (11-13) to the existing cluster
# table in Sheet1. Columns: Rv_ID, annot_int, Gene names, Function [CC],
# cluster_labels. "Function [CC]" has no data for these genes yet, so it
# is left blank (matches the empty <is/> cell in the target row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Rv0554",  3, "bpoC Rv0554",  13),
    @("Rv2940c", 1, "mas Rv2940c",  13),
    @("Rv2808",  1, "Rv2808",       13)
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Write an empty text value into column D (Function [CC]) so the cell
    # is materialized (rather than skipped) but still reads as blank, then
    # reset its style so no stray text-prefix formatting is left behind.
    $ws.Cells.Item($r, 4).Value = "'"
    $ws.Cells.Item($r, 4).Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $row[3]
}
